$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "f"

$ws.Range("B6").Select()
